# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties"
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the other header cells (bold, centered, bordered)
# by copying the format from an existing header cell (AC1) onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Data rows 2-69: Wins=68, Losses=94, Ties=0 for every row.
$lastRow = 69
$ws.Range("AD2:AD$lastRow").Value = 68
$ws.Range("AE2:AE$lastRow").Value = 94
$ws.Range("AF2:AF$lastRow").Value = 0
